$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Columns A and D hold text that Excel would otherwise auto-convert
# (A18 looks like a date, D18 looks like a plain number "00").
# Forcing a text NumberFormat before assignment keeps them as literal
# strings; ClearFormats afterwards drops the now-unneeded explicit
# style so the cell matches the unstyled cells used elsewhere in the
# sheet for this column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-04"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "21:26:04"
$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 140545
$ws.Cells.Item($row, 6).Value = 142870
$ws.Cells.Item($row, 7).Value = 172312
$ws.Cells.Item($row, 8).Value = 147207
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118003
$ws.Cells.Item($row, 11).Value = 224377
$ws.Cells.Item($row, 12).Value = 248583
$ws.Cells.Item($row, 13).Value = 184634
$ws.Cells.Item($row, 14).Value = 110069
$ws.Cells.Item($row, 15).Value = 40401
$ws.Cells.Item($row, 16).Value = 30803
$ws.Cells.Item($row, 17).Value = 72387
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41766
$ws.Cells.Item($row, 20).Value = -1
